# Auto-generated edit script: apply value updates described by the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1103.8
$ws.Range("I92").Value = 1075.5
$ws.Range("J92").Value = 1217
$ws.Range("K92").Value = 1075.5
$ws.Range("L92").Value = 1217
$ws.Range("M92").Value = 172.5
$ws.Range("N92").Value = -3713

$ws.Range("H113").Value = 2999.5454
$ws.Range("I113").Value = 1600.2
$ws.Range("J113").Value = 4165.6665
$ws.Range("K113").Value = 1600.2
$ws.Range("L113").Value = 4165.6665
$ws.Range("M113").Value = 1653.8
$ws.Range("N113").Value = -10673.6665

$ws.Range("H116").Value = 6558.375
$ws.Range("I116").Value = 8062.9443
$ws.Range("J116").Value = 2044.6666
$ws.Range("K116").Value = 8062.9443
$ws.Range("L116").Value = 2044.6666
$ws.Range("M116").Value = -4620.9443
$ws.Range("N116").Value = -8928.6666

$ws.Range("H134").Value = 41666.668
$ws.Range("J134").Value = 41666.668
$ws.Range("L134").Value = 41666.668
$ws.Range("N134").Value = -51806.668

$ws.Range("H135").Value = 2888.6553
$ws.Range("I135").Value = 1572.2916
$ws.Range("J135").Value = 9207.200000000001
$ws.Range("K135").Value = 14150.6244
$ws.Range("L135").Value = 82864.8
$ws.Range("M135").Value = -11615.6244
$ws.Range("N135").Value = -87934.8

$ws.Range("H137").Value = 1736.8334
$ws.Range("I137").Value = 1144.6
$ws.Range("J137").Value = 2159.8572
$ws.Range("K137").Value = 3433.8
$ws.Range("L137").Value = 6479.571599999999
$ws.Range("M137").Value = -883.7999999999997
$ws.Range("N137").Value = -11579.5716

$ws.Range("H138").Value = 3214.68
$ws.Range("I138").Value = 2223.2727
$ws.Range("J138").Value = 3993.6428
$ws.Range("K138").Value = 6669.8181
$ws.Range("L138").Value = 11980.9284
$ws.Range("M138").Value = -1529.8181
$ws.Range("N138").Value = -22260.9284

$ws.Range("H140").Value = 96715.38
$ws.Range("J140").Value = 96715.38
$ws.Range("L140").Value = 96715.38
$ws.Range("N140").Value = -107075.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 6600
$ws.Range("J44").Value = 6600
$ws.Range("L44").Value = 6600
$ws.Range("N44").Value = -7576

$ws.Range("H61").Value = 2289.6333
$ws.Range("I61").Value = 2270.5356
$ws.Range("J61").Value = 2557
$ws.Range("K61").Value = 2270.5356
$ws.Range("L61").Value = 2557
$ws.Range("M61").Value = -2058.5356
$ws.Range("N61").Value = -2981

$ws.Range("H74").Value = 923.5
$ws.Range("I74").Value = 839.1429000000001
$ws.Range("K74").Value = 839.1429000000001
$ws.Range("M74").Value = 34.85709999999995

$ws.Range("H77").Value = 923.5
$ws.Range("I77").Value = 839.1429000000001
$ws.Range("K77").Value = 4195.7145
$ws.Range("M77").Value = 172.2855

$ws.Range("H97").Value = 1985.95
$ws.Range("I97").Value = 1679.9395
$ws.Range("J97").Value = 3428.5715
$ws.Range("K97").Value = 1679.9395
$ws.Range("L97").Value = 3428.5715
$ws.Range("M97").Value = -1183.9395
$ws.Range("N97").Value = -4420.5715

$ws.Range("H102").Value = 4576.25
$ws.Range("I102").Value = 2122
$ws.Range("J102").Value = 8666.666999999999
$ws.Range("K102").Value = 2122
$ws.Range("L102").Value = 8666.666999999999
$ws.Range("M102").Value = -500
$ws.Range("N102").Value = -11910.667

$ws.Range("H132").Value = 2152.9456
$ws.Range("I132").Value = 1747.7805
$ws.Range("J132").Value = 3339.5
$ws.Range("K132").Value = 5243.3415
$ws.Range("L132").Value = 10018.5
$ws.Range("M132").Value = -2713.3415
$ws.Range("N132").Value = -15078.5

$ws.Range("H136").Value = 2289.6333
$ws.Range("I136").Value = 2270.5356
$ws.Range("J136").Value = 2557
$ws.Range("K136").Value = 6811.6068
$ws.Range("L136").Value = 7671
$ws.Range("M136").Value = -4261.6068
$ws.Range("N136").Value = -12771

$ws.Range("H139").Value = 43857.5
$ws.Range("J139").Value = 43857.5
$ws.Range("L139").Value = 43857.5
$ws.Range("N139").Value = -54137.5

$ws.Range("H140").Value = 108381
$ws.Range("J140").Value = 108381
$ws.Range("L140").Value = 108381
$ws.Range("N140").Value = -118741

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12528.167
$ws.Range("I31").Value = 5291.263
$ws.Range("J31").Value = 15881.854
$ws.Range("K31").Value = 5291.263
$ws.Range("L31").Value = 15881.854
$ws.Range("M31").Value = -4996.263
$ws.Range("N31").Value = -16471.854

$ws.Range("H34").Value = 12528.167
$ws.Range("I34").Value = 5291.263
$ws.Range("J34").Value = 15881.854
$ws.Range("K34").Value = 5291.263
$ws.Range("L34").Value = 15881.854
$ws.Range("M34").Value = -5089.263
$ws.Range("N34").Value = -16285.854

$ws.Range("H58").Value = 1096.6086
$ws.Range("I58").Value = 641.75
$ws.Range("J58").Value = 2136.2856
$ws.Range("K58").Value = 641.75
$ws.Range("L58").Value = 2136.2856
$ws.Range("M58").Value = -438.75
$ws.Range("N58").Value = -2542.2856

$ws.Range("H99").Value = 2561.2222
$ws.Range("I99").Value = 2561.2222
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2561.2222
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1063.2222
$ws.Range("N99").Value = ""

$ws.Range("H126").Value = 2561.2222
$ws.Range("I126").Value = 2561.2222
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7683.6666
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5213.6666
$ws.Range("N126").Value = ""

$ws.Range("H136").Value = 1096.6086
$ws.Range("I136").Value = 641.75
$ws.Range("J136").Value = 2136.2856
$ws.Range("K136").Value = 1925.25
$ws.Range("L136").Value = 6408.8568
$ws.Range("M136").Value = 624.75
$ws.Range("N136").Value = -11508.8568

$ws.Range("H138").Value = 49956
$ws.Range("J138").Value = 49956
$ws.Range("L138").Value = 49956
$ws.Range("N138").Value = -60236

$ws.Range("H140").Value = 65000.5
$ws.Range("J140").Value = 65000.5
$ws.Range("L140").Value = 65000.5
$ws.Range("N140").Value = -75360.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 652.125
$ws.Range("I109").Value = 479.55
$ws.Range("J109").Value = 1515
$ws.Range("K109").Value = 1438.65
$ws.Range("L109").Value = 4545
$ws.Range("M109").Value = -398.6500000000001
$ws.Range("N109").Value = -6625

$ws.Range("H115").Value = 3940.2415
$ws.Range("J115").Value = 4599.7085
$ws.Range("L115").Value = 13799.1255
$ws.Range("N115").Value = -16149.1255

$ws.Range("H131").Value = 7463657
$ws.Range("I131").Value = 1218.6666
$ws.Range("J131").Value = 9616283
$ws.Range("K131").Value = 3655.9998
$ws.Range("L131").Value = 28848849
$ws.Range("M131").Value = 1384.0002
$ws.Range("N131").Value = -28858929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3676.2
$ws.Range("I97").Value = 1882.5
$ws.Range("J97").Value = 4872
$ws.Range("K97").Value = 1882.5
$ws.Range("L97").Value = 4872
$ws.Range("M97").Value = -1386.5
$ws.Range("N97").Value = -5864

$ws.Range("H102").Value = 2258.8667
$ws.Range("I102").Value = 2258.8667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2258.8667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -636.8667
$ws.Range("N102").Value = ""

$ws.Range("H113").Value = 4401621.5
$ws.Range("I113").Value = 7693544.5
$ws.Range("J113").Value = 835371.8
$ws.Range("K113").Value = 7693544.5
$ws.Range("L113").Value = 835371.8
$ws.Range("M113").Value = -7691374.5
$ws.Range("N113").Value = -839711.8

$ws.Range("H126").Value = 18520220
$ws.Range("I126").Value = 1600
$ws.Range("K126").Value = 4800
$ws.Range("M126").Value = -2330

$ws.Range("H132").Value = 3026
$ws.Range("I132").Value = 2107.5833
$ws.Range("J132").Value = 4862.8335
$ws.Range("K132").Value = 6322.749899999999
$ws.Range("L132").Value = 14588.5005
$ws.Range("M132").Value = -3792.749899999999
$ws.Range("N132").Value = -19648.5005

$ws.Range("H138").Value = 69450
$ws.Range("J138").Value = 69450
$ws.Range("L138").Value = 69450
$ws.Range("N138").Value = -79730

$ws.Range("H139").Value = 52564.832
$ws.Range("J139").Value = 52564.832
$ws.Range("L139").Value = 52564.832
$ws.Range("N139").Value = -62844.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4161.25
$ws.Range("I7").Value = 3066.75
$ws.Range("J7").Value = 5255.75
$ws.Range("K7").Value = 3066.75
$ws.Range("L7").Value = 5255.75
$ws.Range("M7").Value = -2954.75
$ws.Range("N7").Value = -5479.75

$ws.Range("H93").Value = 4653.1816
$ws.Range("I93").Value = 4057.5715
$ws.Range("J93").Value = 5695.5
$ws.Range("K93").Value = 4057.5715
$ws.Range("L93").Value = 5695.5
$ws.Range("M93").Value = -2809.5715
$ws.Range("N93").Value = -8191.5

$ws.Range("H126").Value = 4161.25
$ws.Range("I126").Value = 3066.75
$ws.Range("J126").Value = 5255.75
$ws.Range("K126").Value = 9200.25
$ws.Range("L126").Value = 15767.25
$ws.Range("M126").Value = -6730.25
$ws.Range("N126").Value = -20707.25

$ws.Range("H136").Value = 4671.3706
$ws.Range("I136").Value = 3609.3333
$ws.Range("J136").Value = 5521
$ws.Range("K136").Value = 10827.9999
$ws.Range("L136").Value = 16563
$ws.Range("M136").Value = -8277.999899999999
$ws.Range("N136").Value = -21663

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2499
$ws.Range("I126").Value = 2164.8333
$ws.Range("J126").Value = 3501.5
$ws.Range("K126").Value = 6494.499899999999
$ws.Range("L126").Value = 10504.5
$ws.Range("M126").Value = -4024.499899999999
$ws.Range("N126").Value = -15444.5

$ws.Range("H136").Value = 986.6389
$ws.Range("I136").Value = 638.2857
$ws.Range("J136").Value = 1474.3334
$ws.Range("K136").Value = 1914.8571
$ws.Range("L136").Value = 4423.0002
$ws.Range("M136").Value = 635.1428999999998
$ws.Range("N136").Value = -9523.0002

$ws.Range("H138").Value = 75016.664

$ws.Range("H139").Value = 53376.875
$ws.Range("J139").Value = 53376.875
$ws.Range("L139").Value = 53376.875
$ws.Range("N139").Value = -63656.875

$ws.Range("H140").Value = 59800
$ws.Range("J140").Value = 59800
$ws.Range("L140").Value = 59800
$ws.Range("N140").Value = -70160
